$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the price/volume columns as text first so the updated values
# are stored as literal strings (matching the original inline-string cells)
# rather than being auto-converted to numbers/percentages.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value2 = "307.98"
$ws.Range("E2").Value2 = "-0.75%"
$ws.Range("D3").Value2 = "37.00"
$ws.Range("E3").Value2 = "-0.38%"
$ws.Range("D4").Value2 = "5.112"
$ws.Range("E4").Value2 = "-0.22%"
$ws.Range("D5").Value2 = "0.07823"
$ws.Range("E5").Value2 = "0.83%"
$ws.Range("D6").Value2 = "8.262"
$ws.Range("E6").Value2 = "0.65%"
$ws.Range("D7").Value2 = "1.882"
$ws.Range("E7").Value2 = "0.21%"
$ws.Range("D8").Value2 = "2.988"
$ws.Range("E8").Value2 = "1.94%"
$ws.Range("D9").Value2 = "0.9215"
$ws.Range("E9").Value2 = "0.26%"
$ws.Range("D10").Value2 = "0.1077"
$ws.Range("E10").Value2 = "-11.37%"
$ws.Range("D11").Value2 = "0.1899"
$ws.Range("E11").Value2 = "0.17%"
$ws.Range("D12").Value2 = "0.08916"
$ws.Range("E12").Value2 = "-2.66%"
$ws.Range("D13").Value2 = "0.03323"
$ws.Range("E13").Value2 = "-3.24%"
$ws.Range("D14").Value2 = "0.09575"
$ws.Range("E14").Value2 = "-1.22%"
$ws.Range("D15").Value2 = "0.001382"
$ws.Range("E15").Value2 = "0.96%"
$ws.Range("D16").Value2 = "0.005878"
$ws.Range("E16").Value2 = "-0.46%"
$ws.Range("D17").Value2 = "3.458"
$ws.Range("E17").Value2 = "-2.86%"
$ws.Range("E18").Value2 = "0.27%"
$ws.Range("D19").Value2 = "0.3422"
$ws.Range("E19").Value2 = "0.46%"
$ws.Range("D20").Value2 = "6.279"
$ws.Range("E20").Value2 = "18.72%"
$ws.Range("D21").Value2 = "0.1281"
$ws.Range("E21").Value2 = "0.96%"
$ws.Range("D22").Value2 = "0.2454"
$ws.Range("E22").Value2 = "-5.41%"
$ws.Range("D23").Value2 = "0.04319"
$ws.Range("E23").Value2 = "-0.93%"
$ws.Range("D24").Value2 = "0.001195"
$ws.Range("E24").Value2 = "-0.38%"
$ws.Range("D25").Value2 = "0.004244"
$ws.Range("E25").Value2 = "0.01%"
$ws.Range("D26").Value2 = "0.0001301"
$ws.Range("E26").Value2 = "-0.07%"
$ws.Range("E27").Value2 = "-98.11%"
$ws.Range("D39").Value2 = "0.02145"
$ws.Range("E39").Value2 = "3.72%"
$ws.Range("D40").Value2 = "0.04995"
$ws.Range("E40").Value2 = "-0.71%"
$ws.Range("D41").Value2 = "0.007524"
$ws.Range("E41").Value2 = "-2.05%"
$ws.Range("E42").Value2 = "0.29%"
$ws.Range("D43").Value2 = "0.008650"
$ws.Range("E43").Value2 = "-11.80%"
$ws.Range("D44").Value2 = "0.002042"
$ws.Range("E44").Value2 = "-2.02%"
$ws.Range("D45").Value2 = "0.008795"
$ws.Range("E45").Value2 = "-8.37%"
$ws.Range("D46").Value2 = "0.00006550"
$ws.Range("E46").Value2 = "-2.45%"
$ws.Range("D47").Value2 = "0.00000000751"
$ws.Range("E47").Value2 = "-0.07%"
$ws.Range("D48").Value2 = "0.002827"
$ws.Range("E48").Value2 = "-3.80%"
$ws.Range("E49").Value2 = "-16.62%"
$ws.Range("D50").Value2 = "0.00002102"
$ws.Range("E50").Value2 = "-0.07%"
$ws.Range("D51").Value2 = "0.0002002"
$ws.Range("E51").Value2 = "-0.07%"

# Restore the default (unstyled) cell style now that the text is set,
# so no stray number-format style is left attached to these cells.
$ws.Range("D2:E51").Style = "Normal"
